$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.400.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7189"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07937"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3151"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.91"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08150"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.894.34"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.240"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7110"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.398"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008415"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.399.04"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.127.05"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.746"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1589"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.092"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.43"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.83"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.424"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.291"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.227"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05325"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.941"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7545"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.181"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.276.30"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.762"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.466"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.69"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9059"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.023.90"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.803"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5205"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.516"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4363"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.57%  "

# Row 45/46 content swap (BabyDogeCoin <-> PaxDollar)
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000131"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.49%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.16%  "
